$p = $ppt.ActivePresentation
try {
    $p.ApplyTheme()
} catch {
    Write-Host "ERR ApplyTheme():" $_.Exception.Message
}
try {
    $p.ApplyTheme("")
} catch {
    Write-Host "ERR ApplyTheme(''):" $_.Exception.Message
}
try {
    $p.ApplyTheme("C:\nonexistent\theme1.xml")
} catch {
    Write-Host "ERR ApplyTheme(path):" $_.Exception.Message
}
try {
    $p.ApplyTemplate("C:\nonexistent\template.potx")
} catch {
    Write-Host "ERR ApplyTemplate(path):" $_.Exception.Message
}
